$d = $word.ActiveDocument

# The two logo pictures that live in the headers/footers (BTec_Logo-Orange
# in the headers, PearsonLogo in the footers) had their cosmetic
# wp:docPr/name and pic:cNvPr/name attributes swapped:
#   footers:  image1.png -> image2.png
#   headers:  image2.jpg -> image1.jpg
# InlineShape has no scriptable .Name property (same as real Word), so the
# rename is applied by round-tripping the package's flat WordOpenXML,
# patching the literal name="..." attribute values, and writing it back.

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.WordOpenXML = $xml
